$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_D12 = @'
Background
id="Par1">Internet analytics are increasingly being integrated into public health regulation.

 One specific application is to monitor compliance of website and social media activity with respect to jurisdictional regulations.

 These data may then identify breaches of compliance and inform disciplinary actions.

 Our study aimed to evaluate the novel use of internet analytics by a Canadian chiropractic regulator to determine their registrants compliance with three regulations related to specific health conditions, pregnancy conditions and most recently, claims of improved immunity during the COVID-19 crisis.


Methods
id="Par2">A customized internet search tool (Market Review Tool, MRT) was used by the College of Chiropractors of British Columbia (CCBC), Canada to audit registrants websites and social media activity.

 The audits extracted words whose use within specific contexts is not permitted under CCBC guidelines.

 The MRT was first used in October of 2018 to identify words related to specific health conditions.

 The MRT was again used in December 2019 for words related to pregnancy and most recently in March 2020 for words related to COVID-19. In these three MRT applications, potential cases of word misuse were evaluated by the regulator who then notified the practitioner to comply with existing regulations by a specific date.

 The MRT was then used on that date to determine compliance.

 Those found to be non-compliant were referred to the regulator’s inquiry committee.

 We mapped this process and reported the outcomes with permission of the regulator.


Results
id="Par3">In September 2018, 250 inappropriate mentions of specific health conditions were detected from approximately 1250 registrants with 2 failing to comply.

 The second scan for pregnancy related terms of approximately1350 practitioners revealed 83 inappropriate mentions.

 Following notification, all 83 cases were compliant within the specified timeframe.

 Regarding COVID-19 related words, 97 inappropriate mentions of the word “immune” were detected from 1350 registrants with 7 cases of non-compliance.


Conclusion
id="Par4">Internet analytics are an effective way for regulators to monitor internet activity to protect the public from misleading statements.

 The processes described were effective at bringing about rapid practitioner compliance.

 Given the increasing volume of internet activity by healthcare professionals, internet analytics are an important addition for health care regulators to protect the public they serve.



'@
$ws.Range("D12").Value = $v_D12

$v_E12 = @'
[Greg%Kawchuk%greg.kawchuk@ualberta.ca%1,  Jan%Hartvigsen%jhartvigsen@health.sdu.dk%2,  Jan%Hartvigsen%jhartvigsen@health.sdu.dk%0,  Stan%Innes%S.Innes@murdoch.edu.au%1,  J. Keith%Simpson%k.simpson@murdoch.edu.au%1,  Brian%Gushaty%bgushaty@gushaty.com%1]
'@
$ws.Range("E12").Value = $v_E12

$v_D16 = @'
id="Par1">The aim of this study was to elicit the views of medical faculty students regarding the COVID-19 pandemic.

 This descriptive study was performed with Ondokuz Mayıs University Medical Faculty students on 24–27 March, 2020. The Medical Faculty currently has 2051 students.

 A questionnaire was used as a data collection tool.

 For that purpose, the authors designed a questionnaire specifically for this research via the “Google Forms” web.

 This consisted of 40 open- and close-ended questions.

 The questionnaire was completed by 1375 (67.1%) students.

 Accordingly, 52.4% of medical students reported feeling mentally unwell.

 Although 50.8% of medical students reported generally/usually obtaining information about COVID-19 through the social media, 82.0% did not trust information/messages arriving through the social media and WhatsApp.

 We found that 86.7% of students regarded frequent hand washing as the most important means of protection against COVID-19, and 19.3% of students did not regard the COVID-19 pandemic as a severe public health problem for Turkey at that moment.

 In addition, 61.6% of students stated that a suppression strategy involving tight restrictions need to be applied to being the pandemic under control in Turkey.

 Use can be made of medical students in the transmission of accurate information during the COVID-19 pandemic.

 Students can be excellent activists on these subjects in countries in which medical education is suspended.

 Measures therefore need to be taken concerning the transmission of up to date and accurate information to medical students.


'@
$ws.Range("D16").Value = $v_D16

$v_E16 = @'
[Servet%Aker%servetaker@gmail.com%1,  Özlem%Mıdık%NULL%2,  Özlem%Mıdık%NULL%0]
'@
$ws.Range("E16").Value = $v_E16

$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(16).AutoFit()
